$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 123, pushing the existing data
# (old rows 123-181) down to rows 124-182.
$ws.Rows("123").Insert()

# Populate the newly inserted row 123 with the new record.
$ws.Range("A123").Value = 1
$ws.Range("B123").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C123").Value = "Arica y Parinacota"
$ws.Range("D123").Value = 44574
$ws.Range("E123").Value = 15
$ws.Range("F123").Value = "Fruta"
$ws.Range("G123").Value = 100108
$ws.Range("H123").Value = "Tropicales y subtropicales"
$ws.Range("I123").Value = 100108006
$ws.Range("J123").Value = "Plátano"
$ws.Range("K123").Value = "Barraganete"
$ws.Range("L123").Value = "Primera"
$ws.Range("M123").Value = 120
$ws.Range("N123").Value = 26000
$ws.Range("O123").Value = 27000
$ws.Range("P123").Value = 26500
$ws.Range("Q123").Value = "$/caja 20 kilos"
$ws.Range("R123").Value = "Ecuador"
$ws.Range("S123").Value = 1325
$ws.Range("T123").Value = 20
